# Refactored checks for information.
# Appends six new rows (33-38) to Sheet1 of the Case Data workbook,
# duplicating the DUS UCM / TAIL LIGHTS-REAR LICENSE PLATE case pair
# (case 21TRD09386 / defendant Bunner) with varying plea/finding
# combinations (mirrors the existing repeated-case-pair pattern already
# present for 21TRD09437 earlier in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$startRow = 33
$endRow = 38

# Columns D (statute) and I (fine amount) hold numeric-looking text
# ("4510.111", "4513.05", "0", ...) that must stay text, not be coerced
# to numbers -- format as Text before writing into them.
$ws.Range("D" + $startRow + ":D" + $endRow).NumberFormat = "@"
$ws.Range("I" + $startRow + ":I" + $endRow).NumberFormat = "@"

# New data rows to append starting at row 33
$newRows = @(
    @{ A = "21TRD09386"; B = "Bunner"; C = "DUS UCM"; D = "4510.111"; E = "UCM"; F = "No Contest"; G = "Guilty"; H = 0; I = "0" },
    @{ A = "21TRD09386"; B = "Bunner"; C = "TAIL LIGHTS-REAR LICENSE PLATE"; D = "4513.05"; E = "MM"; F = "No Contest"; G = "Guilty"; H = 0; I = "0" },
    @{ A = "21TRD09386"; B = "Bunner"; C = "DUS UCM"; D = "4510.111"; E = "UCM"; F = $null; G = $null; H = 0; I = "0" },
    @{ A = "21TRD09386"; B = "Bunner"; C = "TAIL LIGHTS-REAR LICENSE PLATE"; D = "4513.05"; E = "MM"; F = $null; G = $null; H = 0; I = "0" },
    @{ A = "21TRD09386"; B = "Bunner"; C = "DUS UCM"; D = "4510.111"; E = "UCM"; F = "Guilty"; G = "Guilty"; H = 0; I = "0" },
    @{ A = "21TRD09386"; B = "Bunner"; C = "TAIL LIGHTS-REAR LICENSE PLATE"; D = "4513.05"; E = "MM"; F = "Guilty"; G = "Guilty"; H = 0; I = "0" }
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    if ($row.F -ne $null) {
        $ws.Cells.Item($r, 6).Value = $row.F
    }
    if ($row.G -ne $null) {
        $ws.Cells.Item($r, 7).Value = $row.G
    }
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
}
